# Update NATMI ligand-receptor edge table with new TPM-based results.
# Replaces the 4 existing data rows (Sending/Target cluster pairs limited to
# FAPs and MuSCs) with 8 data rows that also include the ECs and
# Resolving-Mac sending clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the various expression / specificity metrics (numeric).
$rows = @(
    @("ECs",           "Ngf", "Ntrk1", "FAPs",  1, 0.3333333333333333, 0.27169,             0.81507,   0.04814835840308388, 0.04814835840308387, 3, 1,                  0.08849299999999999, 0.265479, 0.9052313210944106, 0.9052313210944106, 0.02404266317,        0.21638396853,       0.04358540208575079,  0.04358540208575078),
    @("ECs",           "Ngf", "Ntrk1", "MuSCs", 1, 0.3333333333333333, 0.27169,             0.81507,   0.04814835840308388, 0.04814835840308387, 1, 0.3333333333333333, 0.009264333333333334, 0.027793, 0.09476867890558938, 0.09476867890558936, 0.002517026723333334, 0.02265324051,       0.004562956317333092, 0.004562956317333091),
    @("FAPs",          "Ngf", "Ntrk1", "FAPs",  3, 1,                  0.7246536666666668, 2.173961,  0.1284216734542145,  0.1284216734542145,  3, 1,                  0.08849299999999999, 0.265479, 0.9052313210944106, 0.9052313210944106, 0.06412677692433333, 0.577140992319,      0.1162513211181136,   0.1162513211181136),
    @("FAPs",          "Ngf", "Ntrk1", "MuSCs", 3, 1,                  0.7246536666666668, 2.173961,  0.1284216734542145,  0.1284216734542145,  1, 0.3333333333333333, 0.009264333333333334, 0.027793, 0.09476867890558938, 0.09476867890558936, 0.006713433119222224, 0.06042089807300001, 0.01217035233610091,  0.01217035233610091),
    @("MuSCs",         "Ngf", "Ntrk1", "FAPs",  3, 1,                  4.641773333333333,  13.92532,  0.8226057862976579,  0.8226057862976578,  3, 1,                  0.08849299999999999, 0.265479, 0.9052313210944106, 0.9052313210944106, 0.4107644475866666,  3.696880028279999,   0.7446485226701353,   0.7446485226701351),
    @("MuSCs",         "Ngf", "Ntrk1", "MuSCs", 3, 1,                  4.641773333333333,  13.92532,  0.8226057862976579,  0.8226057862976578,  1, 0.3333333333333333, 0.009264333333333334, 0.027793, 0.09476867890558938, 0.09476867890558936, 0.04300293541777778, 0.38702641876,       0.07795726362752262,  0.07795726362752259),
    @("Resolving-Mac", "Ngf", "Ntrk1", "FAPs",  1, 0.3333333333333333, 0.004650666666666667, 0.013952, 0.0008241818450437709, 0.0008241818450437709, 3, 1,            0.08849299999999999, 0.265479, 0.9052313210944106, 0.9052313210944106, 0.0004115514453333333, 0.003703963008,    0.0007460752204110015, 0.0007460752204110015),
    @("Resolving-Mac", "Ngf", "Ntrk1", "MuSCs", 1, 0.3333333333333333, 0.004650666666666667, 0.013952, 0.0008241818450437709, 0.0008241818450437709, 1, 0.3333333333333333, 0.009264333333333334, 0.027793, 0.09476867890558938, 0.09476867890558936, 0.00004308532622222223, 0.0003877679360000001, 0.00007810662463276934, 0.00007810662463276932)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

$wb.Save()
